$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing existing data down
$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "price"
$ws.Range("D1").Value = "category"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "description"

# Append new data rows (rows 5-8)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "kiwi"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "freshproduct"
$ws.Range("E5").Value = "kiwi1"
$ws.Range("F5").Value = "sfjskls"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "cake"
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = "bakery"
$ws.Range("E6").Value = "chocolate"
$ws.Range("F6").Value = "fsdklffklssf"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "cupcake"
$ws.Range("C7").Value = 44
$ws.Range("D7").Value = "bakery"
$ws.Range("E7").Value = "cake"
$ws.Range("F7").Value = "fdsjklsdfd;lk"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "milk"
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = "grocery"
$ws.Range("E8").Value = "milk"
$ws.Range("F8").Value = "fjsdklsdjfs"

[void]$ws.Range("F8").Select()

# Resize columns to fit the new content (matches widths produced by Excel's
# real AutoFit on the final data for each column)
$ws.Columns.Item(1).ColumnWidth = 1.6666666666666667
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 4.166666666666667
$ws.Columns.Item(5).ColumnWidth = 8.333333333333334
$ws.Columns.Item(6).ColumnWidth = 36.666666666666664

$ws.PageSetup.Orientation = 1
